$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new rows 37-45 (level 0 tutorial strings).
# Values are assigned in the same order they were added to the shared
# string table so new entries line up with the target sharedStrings index.

# Row 37
$ws.Range("A37").Value = "intro_0"
$ws.Range("B37").Value = "Welcome to Kero Builder!"

# Row 38 / Row 39 keys first
$ws.Range("A38").Value = "intro_1"
$ws.Range("A39").Value = "intro_2"

# Row 38 / Row 39 values
$ws.Range("B38").Value = "You will be playing as a builder to help these spacefaring frogs settle into their new planet. "
$ws.Range("B39").Value = "Let's do our best to build their homes!"

# Row 41
$ws.Range("A41").Value = "level_0_intro_0_1"
$ws.Range("B41").Value = "Before we begin, let's first look at the view controls."

# Row 42
$ws.Range("A42").Value = "level_0_intro_0_2"
$ws.Range("B42").Value = "You can drag the view around to get a better look at the map."
$ws.Range("C42").Value = 3

# Row 43 key
$ws.Range("A43").Value = "level_0_intro_0_3"

# Row 40 (filled after row 43's key)
$ws.Range("A40").Value = "next_instruct"
$ws.Range("B40").Value = "Press this button when you're ready to continue."
$ws.Range("C40").Value = 3

# Row 43 value
$ws.Range("B43").Value = "These buttons allow you to rotate or elevate the view."
$ws.Range("C43").Value = 3

# Row 44
$ws.Range("A44").Value = "level_0_intro_1_1"
$ws.Range("B44").Value = "The information on the upper-left tells you the measurement of the unit cube."

# Row 45
$ws.Range("A45").Value = "level_0_intro_1_2"

# Match the saved selection/active cell from the edit.
$ws.Range("B45").Select()
